$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (apex, kings) entirely.
$ws.Rows("5:6").Delete()

# Insert a new column H for "company_storage_path", pushing the existing
# H (created_time) and I (last_updated_time) columns to I and J.
$ws.Columns("H:H").Insert()

# New header cell.
$ws.Range("H1").Value = "company_storage_path"

# New per-company storage path values.
$ws.Range("H2").Value = "C:\Users\jovan\Downloads\grp_quotation_generator\server\Final_Doc/GRPT"
$ws.Range("H3").Value = "C:\Users\jovan\Downloads\grp_quotation_generator\server\Final_Doc/GRPPT"
$ws.Range("H4").Value = "C:\Users\jovan\Downloads\grp_quotation_generator\server\Final_Doc/CLX"

# Set the width of the new column (chosen so the exported <col> width
# rounds to 78, matching the saved workbook).
$ws.Columns("H:H").ColumnWidth = 77.17

# Update the selection to match the new cursor position.
$ws.Range("H11").Select()
